$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "estimate" flag for dCit (row 5) from "yes" to "no"
$ws.Range("F5").Value = "no"

# Update p0 for nMperUnit (row 11) and its "estimate" flag from "yes" to "no"
$ws.Range("E11").Value = 12.2533526453471
$ws.Range("F11").Value = "no"

# Update the active cell selection to match the saved view state
$ws.Activate()
$ws.Range("J9").Select()
